# Fixed a bug in Mask
# Rows 2-15 and rows 19-21 of Sheet1 are re-ordered (the underlying data set
# is the same, just the row order/mask differs). Apply the new row order by
# writing the correct values into each row/column explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1001, 18, 30, 75, 60, 72)
    3  = @(201, 9, 30, 15, 45, 30)
    4  = @(1202, 2, 10, 10, 10, 10)
    5  = @(901, 16, 15, 45, 60, 60)
    6  = @(902, 1, 0, 0, 0, 0)
    7  = @(501, 9, 52, 30, 75, 45)
    8  = @(401, 9, 48, 67, 75, 45)
    9  = @(701, 3, 90, 45, 97, 15)
    10 = @(101, 9, 30, 15, 60, 15)
    11 = @(301, 6, 45, 30, 60, 45)
    12 = @(601, 9, 60, 67, 60, 42)
    13 = @(801, 3, 67, 65, 52, 45)
    14 = @(1201, 2, 10, 10, 10, 10)
    15 = @(1203, 3, 15, 15, 15, 15)
    19 = @(502, 0, 4, 0, 0, 0)
    20 = @(802, 0, 4, 5, 4, 0)
    21 = @(3, 0, 3, 3, 3, 3)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
